$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6 (write the "no" label first so shared-string indices line up)
$ws.Range("A6").Value = "no"
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 0

# New column header
$ws.Range("D1").Value = "Getting Married"

# Values for existing rows (2-5) and the new row 6
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 0

# Move selection to E3 as in the final workbook state
$ws.Range("E3").Select()
